$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B and C: set all values (rows 2-9) to 0
$ws.Range("B2:C9").Value = 0

# Column D: set specific new values (rows 2-9)
$ws.Range("D2").Value = 0.680241283951419
$ws.Range("D3").Value = -0.7307375842165882
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("D6").Value = -0.6754520409216147
$ws.Range("D7").Value = 0.7752585287733935
$ws.Range("D8").Value = 0
$ws.Range("D9").Value = 0
